$rowsData = @(
  @("Bathurst", 19153436, 3624702, 1013496, 1164400, 472498, 0, 500000, 0, 25928532),
  @("Campbellton", 10686193, 1889297, 156231, 1093136, 503816, 46500, 405200, 66867, 14847240),
  @("Dieppe", 53323334, 1095459, 345000, 1697000, 978500, 0, 926441, 583007, 58948741),
  @("Edmundston", 23720297, 5564820, 1718846, 1342789, 725500, 5000, 2420632, 12801, 35510685),
  @("Fredericton", 108183452, 1890636, 3336248, 7046293, 1824740, 107000, 1742, 1988495, 124378606),
  @("Miramichi", 26819207, 6268228, 1019512, 804787, 430500, 0, 637309, 19712, 35999255),
  @("Moncton", 143724760, 3108999, 3175793, 10280883, 3110601, 0, 17, 171570, 163572623),
  @("Saint John", 125842974, 17204756, 1620000, 3990746, 5845070, 9493578, 2711300, 0, 166708424),
  @("Beresford", 4498688, 908618, 98324, 54750, 33850, 27400, 0, 54370, 5676000),
  @("Bouctouche", 2955852, 55271, 168453, 663000, 104975, 0, 22400, 129082, 4099033),
  @("Caraquet", 5909755, 1032625, 187021, 120846, 197753, 6000, 0, 0, 7454000),
  @("Dalhousie", 4003564, 1402560, 73573, 230162, 47382, 0, 50000, 34805, 5842046),
  @("Florenceville-Bristol", 2883408, 125423, 213826, 222457, 51548, 4500, 23000, 10302, 3534464),
  @("Grand Bay-Westfield", 5240616, 645014, 314873, 46100, 48880, 0, 0, 55987, 6351470),
  @("Grand Falls/Grand-Sault", 8315976, 793610, 777545, 1000185, 127310, 0, 0, 13312, 11027938),
  @("Hampton", 4798261, 217469, 412546, 436988, 312240, 0, 0, 0, 6177504),
  @("Hartland", 1303484, 43372, 212020, 82350, 37400, 0, 40164, 41850, 1760640),
  @("Lamèque", 1529342, 176117, 234554, 33050, 49800, 0, 0, 1107, 2023970),
  @("Nackawic", 1376477, 34629, 95190, 79500, 6750, 0, 0, 13546, 1606092),
  @("Oromocto", 14925428, 574123, 1312006, 773555, 857057, 15600, 301044, 36263, 18795076),
  @("Quispamsis", 24180623, 114304, 47064, 898315, 445702, 28800, 20947, 48939, 25784694),
  @("Richibucto", 1612895, 87906, 139580, 758084, 98135, 0, 15364, 0, 2711964),
  @("Riverview", 24875683, 2683203, 577600, 1199531, 632514, 0, 0, 71495, 30040026),
  @("Rothesay", 16492464, 128335, 60000, 422300, 96151, 26500, 1000000, 99250, 18325000),
  @("Sackville", 10255924, 67548, 192075, 332400, 447000, 27750, 0, 62888, 11385585),
  @("Saint Andrews", 4144989, 64748, 127558, 132835, 529080, 99585, 50000, 863, 5149658),
  @("Saint-Léonard", 1201898, 334202, 161419, 68100, 96296, 4000, 28000, 0, 1893915),
  @("Saint-Quentin", 1975650, 438307, 147546, 78920, 98204, 47155, 279972, 120512, 3186266),
  @("Shediac", 10585231, 591983, 392293, 588000, 212000, 0, 507500, 193228, 13070235),
  @("Shippagan", 3720684, 449201, 166373, 553770, 33300, 0, 285000, 10376, 5218704),
  @("Saint George", 1949410, 128374, 156837, 4000, 29600, 650, 0, 97303, 2366174),
  @("St. Stephen", 5402445, 1608894, 224617, 736230, 220570, 4200, 923011, 280156, 9400123),
  @("Sussex", 5418706, 572827, 471985, 352500, 85750, 0, 376307, 1197, 7279272),
  @("Woodstock", 7472215, 855370, 641965, 958400, 158300, 0, 114740, 0, 10200990),
  @("Alma", 464367, 4917, 20020, 2000, 46100, 0, 0, 35308, 572712),
  @("Aroostook", 123052, 105851, 0, 200, 0, 0, 0, 30195, 259298),
  @("Atholville", 3332804, 290668, 1454, 76000, 164874, 0, 61126, 222, 3927148),
  @("Balmoral", 1460957, 110900, 161131, 100000, 152322, 0, 0, 10010, 1995320),
  @("Bas-Caraquet", 1072552, 164334, 59210, 16500, 15850, 13500, 85000, 14514, 1441460),
  @("Bath", 321895, 102239, 37592, 9750, 40600, 6755, 0, 0, 518831),
  @("Belledune", 4603552, 345656, 36380, 87700, 27700, 4500, 0, 0, 5105488),
  @("Bertrand", 891948, 199976, 143619, 82120, 43700, 0, 1, 23474, 1384838),
  @("Blacks Harbour", 932916, 270935, 147396, 112700, 14000, 8600, 5000, 0, 1491547),
  @("Blackville", 700122, 208772, 212364, 0, 10500, 0, 0, 0, 1131758),
  @("Cambridge-Narrows", 1029235, 752, 10911, 0, 9328, 0, 0, 29898, 1080124),
  @("Canterbury", 195287, 96202, 66066, 0, 0, 0, 0, 798, 358353),
  @("Cap-Pelé", 2817898, 200933, 100706, 429750, 46750, 0, 132600, 0, 3728637),
  @("Centreville", 520618, 34800, 97698, 14000, 4000, 0, 0, 64100, 735216),
  @("Charlo", 1277895, 98132, 138171, 0, 47950, 0, 5177, 17871, 1585196),
  @("Chipman", 963249, 334810, 143147, 24000, 50911, 0, 0, 59217, 1575334),
  @("Doaktown", 881673, 219992, 35642, 67700, 17396, 1500, 25000, 0, 1248903),
  @("Dorchester", 939646, 91371, 65984, 0, 7975, 6000, 43011, 179134, 1333121),
  @("Drummond", 702250, 4417, 150370, 6000, 6036, 0, 0, 1322, 870395),
  @("Eel River Crossing", 1072353, 394022, 87495, 0, 9091, 0, 9299, 10025, 1582285),
  @("Fredericton Junction", 564416, 73409, 209938, 0, 48095, 0, 34814, 0, 930672),
  @("Gagetown", 909482, 2438, 46520, 0, 12400, 0, 50000, 60588, 1081428),
  @("Grand Manan", 2581194, 275402, 0, 46550, 15125, 22000, 0, 46548, 2986819),
  @("Grande-Anse", 552864, 445807, 178898, 5982, 6930, 0, 2900, 2048, 1195429),
  @("Harvey", 274358, 42490, 34857, 0, 13410, 1750, 0, 23090, 389955),
  @("Hillsborough", 1164914, 259194, 94631, 102500, 3200, 0, 114658, 76823, 1815920),
  @("Lac Baker", 720312, 1012, 32000, 0, 30500, 15000, 0, 7989, 806813),
  @("Le Goulet", 430865, 239999, 11577, 300, 24900, 0, 0, 7032, 714673),
  @("Maisonnette", 393333, 70424, 14099, 9700, 3400, 1706, 0, 9507, 502169),
  @("McAdam", 805978, 545795, 15258, 12000, 36500, 0, 0, 0, 1415531),
  @("Meductic", 194359, 22383, 110598, 6600, 12428, 1575, 0, 0, 347943),
  @("Memramcook", 4520689, 1155800, 71211, 841927, 61875, 4888, 221285, 135349, 7013024),
  @("Millville", 144078, 84603, 0, 300, 35150, 0, 0, 34887, 299018),
  @("Minto", 1403168, 976038, 267439, 35000, 44523, 0, 0, 0, 2726168),
  @("Néguac", 1818370, 291915, 183600, 4400, 148147, 0, 44110, 51463, 2542005),
  @("New Maryland", 4598492, 6235, 140302, 112500, 55800, 0, 25000, 59, 4938388),
  @("Nigadoo", 768592, 115381, 35439, 16800, 16050, 0, 17000, 17000, 986262),
  @("Norton", 913161, 237290, 272710, 0, 31600, 0, 10000, 387, 1465148),
  @("Paquetville", 734831, 9015, 129717, 0, 118850, 0, 0, 0, 992413),
  @("Perth-Andover", 1476035, 275489, 231344, 244964, 20680, 0, 110451, 114476, 2473439),
  @("Petit-Rocher", 1626046, 483549, 105616, 142511, 15472, 2824, 7152, 30427, 2413597),
  @("Petitcodiac", 1192350, 339215, 191593, 176000, 51000, 0, 0, 109622, 2059780),
  @("Plaster Rock", 1090358, 314018, 132231, 60600, 23000, 0, 0, 33742, 1653949),
  @("Pointe-Verte", 700961, 164947, 58243, 23280, 7720, 0, 47771, 21000, 1023922),
  @("Port Elgin", 366842, 36715, 157649, 3296, 26460, 3296, 0, 34734, 628992),
  @("Rexton", 1170119, 7880, 192429, 0, 82586, 55000, 0, 0, 1508014),
  @("Riverside-Albert", 324837, 20426, 118825, 15000, 800, 0, 48115, 12123, 540126),
  @("Rivière-Verte", 549720, 95047, 30403, 22000, 3500, 0, 0, 0, 700670),
  @("Rogersville", 863902, 444041, 179059, 0, 48527, 0, 0, 5328, 1540857),
  @("Saint-Antoine", 1488961, 329914, 162486, 25500, 20451, 0, 0, 8028, 2035340),
  @("Saint-Isidore", 715362, 25983, 356347, 0, 16800, 0, 62000, 0, 1176492),
  @("Saint-Léolin", 324604, 225717, 13867, 0, 5100, 0, 0, 3701, 572989),
  @("Saint-Louis-de-Kent", 913882, 108321, 119001, 115730, 2900, 0, 50000, 0, 1309834),
  @("Sainte-Anne-de-Madawaska", 656431, 195186, 99000, 0, 25000, 0, 0, 0, 975617),
  @("Sainte-Marie-Saint-Raphaël", 523695, 240102, 73714, 0, 11176, 0, 0, 17810, 866497),
  @("Salisbury", 1810331, 204940, 122089, 0, 7500, 0, 150000, 65886, 2360746),
  @("St. Martins", 273328, 1977, 186657, 34950, 23700, 0, 3680, 50176, 574468),
  @("Stanley", 440259, 5621, 88889, 0, 51601, 0, 0, 5564, 591934),
  @("Sussex Corner", 1252360, 38578, 11100, 5000, 16475, 0, 79500, 0, 1403013),
  @("Tide Head", 923650, 42383, 0, 0, 500, 0, 0, 38240, 1004773),
  @("Tracy", 307946, 155058, 4786, 0, 8940, 0, 0, 0, 476730)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 2 to make room for Bathurst, Campbellton, Dieppe, Edmundston.
# This shifts old rows 2-92 down to rows 6-96 (values and formatting both move).
$ws.Rows.Item(2).Resize(4).Insert()

# The newly inserted rows 2-5 lack the numeric/text formatting used by the rest of the
# data rows. Copy the formatting (styles) from row 6 (A6:J6), which still has the
# original data-row formatting, down into the 4 new rows.
$ws.Range("A6:J6").Copy()
$ws.Range("A2:J5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write out the full, corrected data set (95 rows) into A2:J96.
for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $i + 2
    $row = $rowsData[$i]
    $ws.Cells.Item($r, 1).Value() = $row[0]
    $ws.Cells.Item($r, 2).Value() = $row[1]
    $ws.Cells.Item($r, 3).Value() = $row[2]
    $ws.Cells.Item($r, 4).Value() = $row[3]
    $ws.Cells.Item($r, 5).Value() = $row[4]
    $ws.Cells.Item($r, 6).Value() = $row[5]
    $ws.Cells.Item($r, 7).Value() = $row[6]
    $ws.Cells.Item($r, 8).Value() = $row[7]
    $ws.Cells.Item($r, 9).Value() = $row[8]
    $ws.Cells.Item($r, 10).Value() = $row[9]
}

# Resize the table (ListObject) and its autofilter to cover the new range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J96"))

Write-Host "Edit complete"
